$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-11-27 01:19:27"

# Insert a new row at position 7 - this pushes the old rows 7..15 down to 8..16
$ws.Rows.Item(7).Insert()

# Refresh the "taken at" timestamp for every data row (2..16) in column A
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# Populate the brand-new row 7 with the new listing
$ws.Range("B7").Value = "GoogleAppSheetで行政書士向け案件管理アプリ開発"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5441932"
$ws.Range("G7").Value = 123
$ws.Range("H7").Value = "◆開発 ◇アプリ"

# The row insert does not carry the hyperlink relationships along with it, so
# rebuild the full F2:F16 hyperlink set (URL text matches the target URL for
# every row in this sheet) and re-apply the Hyperlink style.
$ws.Hyperlinks.Delete()

$urls = @{
    2  = "https://www.lancers.jp/work/detail/5441470"
    3  = "https://www.lancers.jp/work/detail/5217096"
    4  = "https://www.lancers.jp/work/detail/5441252"
    5  = "https://www.lancers.jp/work/detail/5441207"
    6  = "https://www.lancers.jp/work/detail/5441442"
    7  = "https://www.lancers.jp/work/detail/5441932"
    8  = "https://www.lancers.jp/work/detail/5441612"
    9  = "https://www.lancers.jp/work/detail/5441557"
    10 = "https://www.lancers.jp/work/detail/5440957"
    11 = "https://www.lancers.jp/work/detail/5441568"
    12 = "https://www.lancers.jp/work/detail/5441448"
    13 = "https://www.lancers.jp/work/detail/5441440"
    14 = "https://www.lancers.jp/work/detail/5440230"
    15 = "https://www.lancers.jp/work/detail/5441609"
    16 = "https://www.lancers.jp/work/detail/5441422"
}

for ($r = 2; $r -le 16; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $cell.Value = $urls[$r]
    $ws.Hyperlinks.Add($cell, $urls[$r])
    $cell.Style = "Hyperlink"
}
